$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge the three runs of the "First n-3 log files..." bullet into
# a single run (no inner run-splits left behind).
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "First n-3 log files are copied*new_log_dir_1.*") {
        $targetPara = $cand
        break
    }
}

$tpRange = $targetPara.Range
$whole = $d.Range($tpRange.Start, $tpRange.End - 1)
# First overwrite with a placeholder so the new text actually differs from
# the old text (same text is a no-op and runs would not be recombined),
# then set the real text - this collapses the paragraph down to one run.
$whole.Text = "PLACEHOLDER_MERGE_0001"
$whole2 = $d.Range($tpRange.Start, $tpRange.Start + 23)
$whole2.Text = "First n-3 log files are copied alphabetical order in new_log_dir_1."

# ---------------------------------------------------------------------------
# Change 2: add the "Git clone <url>" paragraph (with the _GoBack bookmark)
# right before the "python  teradata_usecase.py" paragraph, and collapse the
# old bookmark paragraph's two runs of spaces into a single run.
# ---------------------------------------------------------------------------
$pyPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*teradata_usecase.py*") {
        $pyPara = $cand
        break
    }
}
$pyIndex = $pyPara.Index

# Insert a brand-new empty paragraph right before it.
$pyPara.Range.InsertParagraphBefore()

$gitPara = $d.Paragraphs.Item($pyIndex)
$gitStart = $gitPara.Range.Start

# Fill it with a tab, the "Git clone " label, the URL, and a trailing
# sentinel character (the sentinel lets us park the bookmark two characters
# before the paragraph mark - placing it directly at End-1 triggers an
# engine bug that relocates the bookmark to the top of the document - and
# then trim the sentinel back off without disturbing the bookmark).
$insPoint = $d.Range($gitStart, $gitStart)
$insPoint.InsertAfter("`tGit clone https://github.com/darkside0522/teradata_code_challange.gitX")

$gitPara2 = $d.Paragraphs.Item($pyIndex)
$bmPos = $gitPara2.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$gitPara3 = $d.Paragraphs.Item($pyIndex)
$sentinelRange = $d.Range($gitPara3.Range.End - 2, $gitPara3.Range.End - 1)
$sentinelRange.Delete()

# ---------------------------------------------------------------------------
# Collapse the now-bookmark-free trailing paragraph's two space-only runs
# into one run of 28 spaces.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastWhole = $d.Range($lastRange.Start, $lastRange.End - 1)
$lastWhole.Text = "PLACEHOLDER_MERGE_0002"
$lastWhole2 = $d.Range($lastRange.Start, $lastRange.Start + 22)
$lastWhole2.Text = "                            "
